$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header values ---
# VALOR MORA (E11): 636000 -> 848000
$ws.Range("E11").Value = 848000

# Cant. Periodos (F13): 3 -> 4
$ws.Range("F13").Value = 4

# --- Insert a new data row (row 19) below the existing last data row (row 18) ---
# This mirrors duplicating the last "Periodo Mora" row and pushes the
# signature block (old rows 23-24) down to rows 24-25.
$ws.Rows("19:19").Insert()

# Duplicate row 18 (values + full formatting) into the newly inserted row 19,
# so row 19 ends up with the same borders/shading that row 18 used to have.
$ws.Range("B18:J18").Copy($ws.Range("B19:J19"))

# Re-style row 18 to match the "interior" row look used by rows 16-17
# (row 18 previously had the "last row" heavier-bottom-border style).
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the "Periodo Mora" values in the data table ---
$ws.Range("E16").Value = "2505"
$ws.Range("E17").Value = "2506"
$ws.Range("E18").Value = "2507"
$ws.Range("E19").Value = "2508"
